$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
